$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep column C as plain text (it already is inline text, not a date)
$ws.Range("C2:C9").NumberFormat = "@"

# Row 2: Bahia -> Pernambuco
$ws.Range("A2").Value = "Pernambuco"
$ws.Range("C2").Value = "01/04/2024"
$ws.Range("D2").Value = 11.5

# Row 3: Pernambuco -> Bahia
$ws.Range("A3").Value = "Bahia"
$ws.Range("C3").Value = "01/04/2024"
$ws.Range("D3").Value = 11.1

# Row 4: Amapá -> Distrito Federal
$ws.Range("A4").Value = "Distrito Federal"
$ws.Range("C4").Value = "01/04/2024"
$ws.Range("D4").Value = 9.699999999999999

# Row 5: Rio de Janeiro (name unchanged)
$ws.Range("C5").Value = "01/04/2024"
$ws.Range("D5").Value = 9.6

# Row 6: Piauí -> Rio Grande do Norte
$ws.Range("A6").Value = "Rio Grande do Norte"
$ws.Range("C6").Value = "01/04/2024"
$ws.Range("D6").Value = 9.1

# Row 7: Sergipe (name unchanged)
$ws.Range("C7").Value = "01/04/2024"
$ws.Range("D7").Value = 9.1

# Row 8: Nordeste (name unchanged)
$ws.Range("C8").Value = "01/04/2024"
$ws.Range("D8").Value = 9.4

# Row 9: Brasil (name unchanged)
$ws.Range("C9").Value = "01/04/2024"
$ws.Range("D9").Value = 6.9
